$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing date-formatted cell style down onto the two new date cells
# first, so the new B10/B11 values pick up the same short-date style (s=1)
# used throughout the rest of the column instead of minting a new number format.
$ws.Range("B9").Copy()
$ws.Range("B10:B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New row 10: further work on USB port numbers
$ws.Range("A10").Value = "Ports USB"
$ws.Range("B10").Value = "2/6/2019"
$ws.Range("D10").Value = "J'ai chercher une solution pour récupérer les n° de ports et de hub, j'ai trouvé une librairie mais je n'ai pas réussi à l'utiliser "
$ws.Range("C10").Value = "4h00"

# Fix typos in existing entries
$ws.Range("A6").Value = "Utilisation des threads"

# New row 11: database connection bug fixes
$ws.Range("C11").Value = "2h30"
$ws.Range("D11").Value = "J'ai corrigé les erreurs qui restaient par rapport à la connexion à la base de données, maintenant les enregistrement sont bien créer quand on insère une nouvelle clé, il ne reste plus qu'à mettre à jour les table si on change la clé de port "
$ws.Range("A11").Value = "MySQL"
$ws.Range("B11").Value = "2/7/2019"

$ws.Range("A7").Value = "MySQL"

# Update view state: scroll back to A1, update selection
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("N11").Select()
